$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $find)
    }
}

# 1. Center the title paragraph (first paragraph)
$d.Paragraphs.Item(1).Range.ParagraphFormat.Alignment = 1

# 2. "проекта" -> "предлагаемого в проекте подхода" (only the specific instance)
Replace-Text "Основополагающая идея проекта заключается" "Основополагающая идея предлагаемого в проекте подхода заключается"

# 3. Simplify the "base notions" sentence
Replace-Text "базируется на следующих понятиях. Узел гетерогенной среды представляется вычислительным устройством с операционной исполнительной средой. " "базируется на понятиях задачи и службы мониторинга. "

# 4. Insert "гетерогенной среды" before period
Replace-Text "состоянии удаленного узла. Служба мониторинга" "состоянии удаленного узла гетерогенной среды. Служба мониторинга"

# 5. Remove storage explanation sentence
Replace-Text " в ассоциированном хранилище данных. Наконец, хранилище данных - это пассивная сущность, предоставляющая службам ресурсы для приема и сохранения сообщений." " в хранилище данных."

# 6. "с использованием ресурсов операционной среды." -> "на базе ресурсов операционной среды." then split paragraph
Replace-Text " с использованием ресурсов операционной среды. " " на базе ресурсов операционной среды.^p "

# 8. "Это состояние" -> "В классической трактовке, это состояние"
Replace-Text "момент (рисунок 2). Это состояние" "момент (рисунок 2). В классической трактовке, это состояние"

# 9. "нагрузкой на узлы" -> "величиной нагрузки на узлы"
Replace-Text "экземпляров модулей и нагрузкой на узлы." "экземпляров модулей и величиной нагрузки на узлы."

# 10. "В предлагаемой архитектуре сущность" -> "В предлагаемой архитектуре системы мониторинга сущность"
Replace-Text "злы. В предлагаемой архитектуре сущность распределенного модуля" "злы. В предлагаемой архитектуре системы мониторинга сущность распределенного модуля"

# 11. "ей некоторые особенности элемента распределенной системы: репликация" -> "ей такие особенности элемента распределенной системы как репликация"
Replace-Text "ей некоторые особенности элемента распределенной системы: репликация" "ей такие особенности элемента распределенной системы как репликация"

# 12. "переносимость. " -> "переносимость." (trailing space trimmed at paragraph end)
Replace-Text "репликация, масштабируемость и переносимость. " "репликация, масштабируемость и переносимость."

# 13. Insert extra responsibility clause
Replace-Text "планирование запусков модулей мониторинга; предоставление промежуточного хранилища" "планирование запусков модулей мониторинга; мониторинг и диспетчеризация процессов исполнения модулей мониторинга; предоставление промежуточного хранилища"

# 14. Append new sentences about states of the service
Replace-Text " – неопределенное, сетевое и автономное. " " – неопределенное, сетевое и автономное. Служба находится в неопределенном состоянии в данный момент если не обладает никакой информацией о системном окружении. Сетевое или автономное состояния службы определяются соответственно доступностью или отсутствием операционной сетевой подсистемы узла, в рамках которого она запущена."

# 15. Rewrite finite automaton sentence
Replace-Text "Функционирование службы мониторинга можно описать в терминах конечных автоматов. " "В некотором смысле, функционирование службы мониторинга можно описать в терминах конечных автоматов. Тогда, служба мониторинга представляет собой детерминированный автомат с конечным набором состояний, а переходы между состояниями осуществляются по наступлению некоторого внутреннего события из множества допустимых событий, генерируемых подсистемами окружения."

# 15b. Merge away the two now-empty BodyText paragraphs that followed it
$autoPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.StartsWith("В некотором смысле")) {
        $autoPara = $d.Paragraphs.Item($i)
        break
    }
}
if ($autoPara -ne $null) {
    $endPos = $autoPara.Range.End
    $d.Range($endPos - 1, $endPos).Delete()
    $d.Range($endPos - 1, $endPos).Delete()
}

# 17. Box set wording
Replace-Text "поддерживаемых прикладным интерфейсом программирования; разработка шаблонных модулей" "поддерживаемых прикладным интерфейсом программирования модулей; разработка коробочного набора шаблонных модулей"

Write-Output "DONE"
